$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Force the cell to keep its literal text content (e.g. "27.168.71",
    # "216.03", "3.00") instead of Excel auto-converting number-looking
    # strings to numeric values (which would drop separators/zeros).
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.168.71"
Set-TextValue $ws.Range("E2") "  -0.03%  "

Set-TextValue $ws.Range("D3") "1.632.11"
Set-TextValue $ws.Range("E3") "  -0.88%  "

Set-TextValue $ws.Range("E4") "  +0.02%  "

Set-TextValue $ws.Range("D5") "216.03"
Set-TextValue $ws.Range("E5") "  -0.54%  "

Set-TextValue $ws.Range("D6") "0.521"
Set-TextValue $ws.Range("E6") "  +1.10%  "

Set-TextValue $ws.Range("E7") "  +0.04%  "

Set-TextValue $ws.Range("D8") "0.255"
Set-TextValue $ws.Range("E8") "  -0.34%  "

Set-TextValue $ws.Range("D9") "0.0626"
Set-TextValue $ws.Range("E9") "  -0.18%  "

Set-TextValue $ws.Range("D10") "20.23"
Set-TextValue $ws.Range("E10") "  +1.48%  "

Set-TextValue $ws.Range("D11") "0.0850"
Set-TextValue $ws.Range("E11") "  +0.07%  "

Set-TextValue $ws.Range("D12") "1.648.51"
Set-TextValue $ws.Range("E12") "  +0.04%  "

Set-TextValue $ws.Range("D13") "4.14"
Set-TextValue $ws.Range("E13") "  -0.11%  "

Set-TextValue $ws.Range("D14") "0.544"
Set-TextValue $ws.Range("E14") "  +0.67%  "

Set-TextValue $ws.Range("D15") "64.98"
Set-TextValue $ws.Range("E15") "  -3.90%  "

Set-TextValue $ws.Range("D16") "27.169.48"
Set-TextValue $ws.Range("E16") "  -0.02%  "

Set-TextValue $ws.Range("E17") "  +0.44%  "

Set-TextValue $ws.Range("D18") "217.61"
Set-TextValue $ws.Range("E18") "  -0.56%  "

Set-TextValue $ws.Range("D20") "6.92"
Set-TextValue $ws.Range("E20") "  +1.37%  "

Set-TextValue $ws.Range("D21") "4.40"
Set-TextValue $ws.Range("E21") "  -0.23%  "

Set-TextValue $ws.Range("D22") "2.43"
Set-TextValue $ws.Range("E22") "  -6.01%  "

Set-TextValue $ws.Range("D23") "9.05"
Set-TextValue $ws.Range("E23") "  -1.81%  "

Set-TextValue $ws.Range("D24") "147.99"
Set-TextValue $ws.Range("E24") "  +0.05%  "

Set-TextValue $ws.Range("E25") "  +0.07%  "

Set-TextValue $ws.Range("D26") "7.31"
Set-TextValue $ws.Range("E26") "  -3.44%  "

Set-TextValue $ws.Range("E27") "  +0.27%  "

Set-TextValue $ws.Range("D28") "15.65"
Set-TextValue $ws.Range("E28") "  -0.61%  "

Set-TextValue $ws.Range("D29") "0.0507"
Set-TextValue $ws.Range("E29") "  -0.41%  "

Set-TextValue $ws.Range("E30") "  -0.37%  "

Set-TextValue $ws.Range("D31") "3.38"
Set-TextValue $ws.Range("E31") "  -0.59%  "

Set-TextValue $ws.Range("D32") "3.00"

Set-TextValue $ws.Range("D33") "1.341.05"
Set-TextValue $ws.Range("E33") "  +6.10%  "

Set-TextValue $ws.Range("D34") "1.57"
Set-TextValue $ws.Range("E34") "  -0.15%  "

Set-TextValue $ws.Range("E35") "  -0.13%  "

Set-TextValue $ws.Range("E36") "  -0.61%  "

Set-TextValue $ws.Range("D37") "0.546"
Set-TextValue $ws.Range("E37") "  -0.37%  "

Set-TextValue $ws.Range("D38") "0.851"
Set-TextValue $ws.Range("E38") "  +0.27%  "

Set-TextValue $ws.Range("E39") "  -0.02%  "

Set-TextValue $ws.Range("E40") "  +1.23%  "

Set-TextValue $ws.Range("D41") "0.804"
Set-TextValue $ws.Range("E41") "  -0.50%  "

Set-TextValue $ws.Range("D42") "64.64"
Set-TextValue $ws.Range("E42") "  +4.32%  "

Set-TextValue $ws.Range("D43") "5.27"
Set-TextValue $ws.Range("E43") "  -3.06%  "

Set-TextValue $ws.Range("D44") "1.772.25"
Set-TextValue $ws.Range("E44") "  -0.93%  "

Set-TextValue $ws.Range("E45") "  -0.74%  "

Set-TextValue $ws.Range("E46") "  +0.10%  "

Set-TextValue $ws.Range("D47") "0.822"
Set-TextValue $ws.Range("E47") "  +22.76%  "

Set-TextValue $ws.Range("B48") "BabyDogeCoin"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D48") "0.0₆0105"
Set-TextValue $ws.Range("E48") "  -2.08%  "

Set-TextValue $ws.Range("B49") "Cronos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.0513"
Set-TextValue $ws.Range("E49") "  -0.02%  "

Set-TextValue $ws.Range("B50") "Algorand"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D50") "0.0993"
Set-TextValue $ws.Range("E50") "  +2.04%  "

Set-TextValue $ws.Range("B51") "EnergySwap"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.58"
Set-TextValue $ws.Range("E51") "  -1.05%  "
